$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-13 (rows 10-13 are newly added) with recomputed NATMI values
# following Dr Hou advice (adds the ECs target-cluster column and updates stats).

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hras"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.980814333333332
$ws.Range("H2").Value = 20.942443
$ws.Range("I2").Value = 0.2573350203399358
$ws.Range("J2").Value = 0.2573350203399358
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.7246186666666666
$ws.Range("N2").Value = 2.173856
$ws.Range("O2").Value = 0.02655111241446272
$ws.Range("P2").Value = 0.02655111241446271
$ws.Range("Q2").Value = 5.058428374467554
$ws.Range("R2").Value = 45.52585537020799
$ws.Range("S2").Value = 0.006832531053223684
$ws.Range("T2").Value = 0.006832531053223683

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hras"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.980814333333332
$ws.Range("H3").Value = 20.942443
$ws.Range("I3").Value = 0.2573350203399358
$ws.Range("J3").Value = 0.2573350203399358
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 17.195945
$ws.Range("N3").Value = 51.587835
$ws.Range("O3").Value = 0.630085160334334
$ws.Range("P3").Value = 0.630085160334334
$ws.Range("Q3").Value = 120.0416993312116
$ws.Range("R3").Value = 1080.375293980905
$ws.Range("S3").Value = 0.1621429775505275
$ws.Range("T3").Value = 0.1621429775505275

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hras"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.980814333333332
$ws.Range("H4").Value = 20.942443
$ws.Range("I4").Value = 0.2573350203399358
$ws.Range("J4").Value = 0.2573350203399358
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 9.370898
$ws.Range("N4").Value = 28.112694
$ws.Range("O4").Value = 0.3433637272512032
$ws.Range("P4").Value = 0.3433637272512032
$ws.Range("Q4").Value = 65.41649907460466
$ws.Range("R4").Value = 588.7484916714419
$ws.Range("S4").Value = 0.08835951173618453
$ws.Range("T4").Value = 0.08835951173618453

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hras"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.375361666666667
$ws.Range("H5").Value = 19.126085
$ws.Range("I5").Value = 0.2350161092714131
$ws.Range("J5").Value = 0.2350161092714131
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.7246186666666666
$ws.Range("N5").Value = 2.173856
$ws.Range("O5").Value = 0.02655111241446272
$ws.Range("P5").Value = 0.02655111241446271
$ws.Range("Q5").Value = 4.619706070417777
$ws.Range("R5").Value = 41.57735463376
$ws.Range("S5").Value = 0.006239939136474944
$ws.Range("T5").Value = 0.006239939136474943

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hras"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.375361666666667
$ws.Range("H6").Value = 19.126085
$ws.Range("I6").Value = 0.2350161092714131
$ws.Range("J6").Value = 0.2350161092714131
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 17.195945
$ws.Range("N6").Value = 51.587835
$ws.Range("O6").Value = 0.630085160334334
$ws.Range("P6").Value = 0.630085160334334
$ws.Range("Q6").Value = 109.6303685751083
$ws.Range("R6").Value = 986.673317175975
$ws.Range("S6").Value = 0.1480801628914297
$ws.Range("T6").Value = 0.1480801628914297

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hras"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.375361666666667
$ws.Range("H7").Value = 19.126085
$ws.Range("I7").Value = 0.2350161092714131
$ws.Range("J7").Value = 0.2350161092714131
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 9.370898
$ws.Range("N7").Value = 28.112694
$ws.Range("O7").Value = 0.3433637272512032
$ws.Range("P7").Value = 0.3433637272512032
$ws.Range("Q7").Value = 59.74286389144333
$ws.Range("R7").Value = 537.68577502299
$ws.Range("S7").Value = 0.08069600724350846
$ws.Range("T7").Value = 0.08069600724350846

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Hras"
$ws.Range("C8").Value = "Agtr1a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.973131
$ws.Range("H8").Value = 17.919393
$ws.Range("I8").Value = 0.2201886075150976
$ws.Range("J8").Value = 0.2201886075150976
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.7246186666666666
$ws.Range("N8").Value = 2.173856
$ws.Range("O8").Value = 0.02655111241446272
$ws.Range("P8").Value = 0.02655111241446271
$ws.Range("Q8").Value = 4.328242221045333
$ws.Range("R8").Value = 38.954179989408
$ws.Range("S8").Value = 0.005846252470517367
$ws.Range("T8").Value = 0.005846252470517366

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Hras"
$ws.Range("C9").Value = "Agtr1a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.973131
$ws.Range("H9").Value = 17.919393
$ws.Range("I9").Value = 0.2201886075150976
$ws.Range("J9").Value = 0.2201886075150976
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 17.195945
$ws.Range("N9").Value = 51.587835
$ws.Range("O9").Value = 0.630085160334334
$ws.Range("P9").Value = 0.630085160334334
$ws.Range("Q9").Value = 102.713632153795
$ws.Range("R9").Value = 924.422689384155
$ws.Range("S9").Value = 0.138737574069944
$ws.Range("T9").Value = 0.138737574069944

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Hras"
$ws.Range("C10").Value = "Agtr1a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.973131
$ws.Range("H10").Value = 17.919393
$ws.Range("I10").Value = 0.2201886075150976
$ws.Range("J10").Value = 0.2201886075150976
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 9.370898
$ws.Range("N10").Value = 28.112694
$ws.Range("O10").Value = 0.3433637272512032
$ws.Range("P10").Value = 0.3433637272512032
$ws.Range("Q10").Value = 55.97360134163799
$ws.Range("R10").Value = 503.762412074742
$ws.Range("S10").Value = 0.0756047809746362
$ws.Range("T10").Value = 0.0756047809746362

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Hras"
$ws.Range("C11").Value = "Agtr1a"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.798031999999999
$ws.Range("H11").Value = 23.394096
$ws.Range("I11").Value = 0.2874602628735535
$ws.Range("J11").Value = 0.2874602628735535
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.7246186666666666
$ws.Range("N11").Value = 2.173856
$ws.Range("O11").Value = 0.02655111241446272
$ws.Range("P11").Value = 0.02655111241446271
$ws.Range("Q11").Value = 5.650599550463999
$ws.Range("R11").Value = 50.85539595417599
$ws.Range("S11").Value = 0.007632389754246723
$ws.Range("T11").Value = 0.007632389754246722

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Hras"
$ws.Range("C12").Value = "Agtr1a"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.798031999999999
$ws.Range("H12").Value = 23.394096
$ws.Range("I12").Value = 0.2874602628735535
$ws.Range("J12").Value = 0.2874602628735535
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 17.195945
$ws.Range("N12").Value = 51.587835
$ws.Range("O12").Value = 0.630085160334334
$ws.Range("P12").Value = 0.630085160334334
$ws.Range("Q12").Value = 134.09452938024
$ws.Range("R12").Value = 1206.85076442216
$ws.Range("S12").Value = 0.1811244458224328
$ws.Range("T12").Value = 0.1811244458224328

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Hras"
$ws.Range("C13").Value = "Agtr1a"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.798031999999999
$ws.Range("H13").Value = 23.394096
$ws.Range("I13").Value = 0.2874602628735535
$ws.Range("J13").Value = 0.2874602628735535
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 9.370898
$ws.Range("N13").Value = 28.112694
$ws.Range("O13").Value = 0.3433637272512032
$ws.Range("P13").Value = 0.3433637272512032
$ws.Range("Q13").Value = 73.07456247273599
$ws.Range("R13").Value = 657.671062254624
$ws.Range("S13").Value = 0.098703427296874
$ws.Range("T13").Value = 0.098703427296874
